$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("G2").Value = 4
$ws.Range("H2").Value = 3.8
$ws.Range("I2").Value = 1.83
$ws.Range("J2").Value = 4
$ws.Range("L2").Value = 2.38
$ws.Range("AA2").Value = 26
$ws.Range("AB2").Value = 26
$ws.Range("AK2").Value = 17
$ws.Range("AN2").Value = 6
$ws.Range("AO2").Value = 19
$ws.Range("AP2").Value = 21
$ws.Range("AX2").Value = 4.33
$ws.Range("AY2").Value = 9.5
$ws.Range("BA2").Value = 29
$ws.Range("BD2").Value = 151

# Row 3
$ws.Range("G3").Value = 2.7
$ws.Range("I3").Value = 2.88
$ws.Range("J3").Value = 3.6
$ws.Range("Q3").Value = 2.75
$ws.Range("R3").Value = 1.44
$ws.Range("S3").Value = 1.62
$ws.Range("T3").Value = 2.2
$ws.Range("X3").Value = 12
$ws.Range("Z3").Value = 29
$ws.Range("AC3").Value = 6
$ws.Range("AI3").Value = 12
$ws.Range("AP3").Value = 34
$ws.Range("AT3").Value = 2.2
$ws.Range("AY3").Value = 17
$ws.Range("BC3").Value = 351

# Row 4
$ws.Range("G4").Value = 1.83
$ws.Range("H4").Value = 3.7
$ws.Range("I4").Value = 4.2
$ws.Range("J4").Value = 2.4
$ws.Range("L4").Value = 4.5
$ws.Range("X4").Value = 9
$ws.Range("AE4").Value = 13
$ws.Range("AL4").Value = 34
$ws.Range("AO4").Value = 9.5
$ws.Range("AQ4").Value = 29
$ws.Range("AU4").Value = 8
$ws.Range("AZ4").Value = 29

# Row 5
$ws.Range("G5").Value = 2.75
$ws.Range("I5").Value = 2.3
$ws.Range("J5").Value = 3.2
$ws.Range("K5").Value = 2.4
$ws.Range("O5").Value = 1.23
$ws.Range("P5").Value = 3.4
$ws.Range("Q5").Value = 1.73
$ws.Range("R5").Value = 2.08
$ws.Range("Y5").Value = 10
$ws.Range("AH5").Value = 10
$ws.Range("AI5").Value = 13
$ws.Range("AK5").Value = 23
$ws.Range("AP5").Value = 19
$ws.Range("AQ5").Value = 41
$ws.Range("AR5").Value = 51
$ws.Range("AS5").Value = 101
$ws.Range("AT5").Value = 3.5
$ws.Range("AU5").Value = 7
$ws.Range("AW5").Value = 301
$ws.Range("AX5").Value = 4.75
$ws.Range("AZ5").Value = 17
$ws.Range("BC5").Value = 101

# Row 6
$ws.Range("G6").Value = 1.62
$ws.Range("I6").Value = 5.25
$ws.Range("J6").Value = 2.25
$ws.Range("L6").Value = 5.5
$ws.Range("M6").Value = 1.05
$ws.Range("N6").Value = 11
$ws.Range("O6").Value = 1.29
$ws.Range("P6").Value = 3.5
$ws.Range("Q6").Value = 1.95
$ws.Range("R6").Value = 1.85
$ws.Range("S6").Value = 1.4
$ws.Range("T6").Value = 2.75
$ws.Range("U6").Value = 1.91
$ws.Range("V6").Value = 1.8
$ws.Range("AN6").Value = 3.6
$ws.Range("AO6").Value = 8.5
$ws.Range("AS6").Value = 151
$ws.Range("AT6").Value = 2.75
$ws.Range("AU6").Value = 8.5
$ws.Range("AV6").Value = 51
$ws.Range("AY6").Value = 29
$ws.Range("AZ6").Value = 34
$ws.Range("BA6").Value = 101
$ws.Range("BB6").Value = 126
$ws.Range("BC6").Value = 251

# Row 7
$ws.Range("G7").Value = 4.5
$ws.Range("H7").Value = 3.8
$ws.Range("I7").Value = 1.7
$ws.Range("J7").Value = 4.5
$ws.Range("L7").Value = 2.25
$ws.Range("U7").Value = 1.53
$ws.Range("V7").Value = 2.38
$ws.Range("Y7").Value = 15
$ws.Range("AD7").Value = 7.5
$ws.Range("AE7").Value = 12
$ws.Range("AF7").Value = 34
$ws.Range("AG7").Value = 101
$ws.Range("AI7").Value = 10
$ws.Range("AK7").Value = 15
$ws.Range("AM7").Value = 19
$ws.Range("AN7").Value = 6.5
$ws.Range("AO7").Value = 21
$ws.Range("AP7").Value = 23
$ws.Range("AQ7").Value = 67
$ws.Range("AR7").Value = 67
$ws.Range("AS7").Value = 126
$ws.Range("AY7").Value = 8.5
$ws.Range("BA7").Value = 26

# Row 8
$ws.Range("I8").Value = 17
$ws.Range("N8").Value = 26
$ws.Range("W8").Value = 12
$ws.Range("Z8").Value = 7.5
$ws.Range("AG8").Value = 301
$ws.Range("AI8").Value = 67
$ws.Range("AK8").Value = 201
$ws.Range("BA8").Value = 301

# Row 9
$ws.Range("I9").Value = 3.6
$ws.Range("J9").Value = 3
$ws.Range("L9").Value = 4.33
$ws.Range("M9").Value = 1.11
$ws.Range("N9").Value = 6.5
$ws.Range("S9").Value = 1.57
$ws.Range("T9").Value = 2.25
$ws.Range("U9").Value = 2.1
$ws.Range("V9").Value = 1.67
$ws.Range("W9").Value = 6
$ws.Range("AC9").Value = 6.5
$ws.Range("AF9").Value = 67
$ws.Range("AH9").Value = 8.5
$ws.Range("AJ9").Value = 13
$ws.Range("AL9").Value = 34
$ws.Range("AM9").Value = 41
$ws.Range("AQ9").Value = 41
$ws.Range("AS9").Value = 251
$ws.Range("AT9").Value = 2.25
$ws.Range("AU9").Value = 9
$ws.Range("AY9").Value = 21
$ws.Range("AZ9").Value = 34

# Row 15
$ws.Range("G15").Value = 4.2
$ws.Range("I15").Value = 1.85
$ws.Range("M15").Value = 1.06
$ws.Range("N15").Value = 10
$ws.Range("O15").Value = 1.33
$ws.Range("P15").Value = 3.25
$ws.Range("Q15").Value = 2.08
$ws.Range("R15").Value = 1.73
$ws.Range("W15").Value = 11
$ws.Range("AC15").Value = 9
$ws.Range("AH15").Value = 6.5
$ws.Range("AJ15").Value = 9

# Row 16
$ws.Range("G16").Value = 2.25
$ws.Range("H16").Value = 3
$ws.Range("I16").Value = 3.6
$ws.Range("J16").Value = 3.1
$ws.Range("W16").Value = 5.5
$ws.Range("X16").Value = 9
$ws.Range("AO16").Value = 13
$ws.Range("AZ16").Value = 41

# Row 17
$ws.Range("G17").Value = 2.38
$ws.Range("H17").Value = 2.75
$ws.Range("I17").Value = 3.6
$ws.Range("J17").Value = 3.25
$ws.Range("K17").Value = 1.8
$ws.Range("L17").Value = 4.5
$ws.Range("O17").Value = 1.67
$ws.Range("P17").Value = 2.1
$ws.Range("Q17").Value = 3.4
$ws.Range("R17").Value = 1.33
$ws.Range("U17").Value = 2.38
$ws.Range("V17").Value = 1.53
$ws.Range("W17").Value = 5.5
$ws.Range("X17").Value = 9.5
$ws.Range("Z17").Value = 23
$ws.Range("AC17").Value = 5
$ws.Range("AF17").Value = 101
$ws.Range("AJ17").Value = 15
$ws.Range("AK17").Value = 41
$ws.Range("AL17").Value = 41
$ws.Range("AN17").Value = 4
$ws.Range("AO17").Value = 15
$ws.Range("AV17").Value = 101
$ws.Range("AY17").Value = 23
$ws.Range("BA17").Value = 81
$ws.Range("BB17").Value = 151
